$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95 (pushes existing rows 95:110 down to 96:111)
# so a new weekly price observation can be recorded at the top of this block.
$ws.Range("A95").EntireRow.Insert()

$ws.Cells.Item(95, 1).Value = 11
$ws.Cells.Item(95, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(95, 3).Value = "Bíobío"
$ws.Cells.Item(95, 4).Value = 44637
$ws.Cells.Item(95, 5).Value = 8
$ws.Cells.Item(95, 6).Value = 100112043
$ws.Cells.Item(95, 7).Value = "Pepino ensalada"
$ws.Cells.Item(95, 8).Value = "Sin especificar"
$ws.Cells.Item(95, 9).Value = "Primera"
$ws.Cells.Item(95, 10).Value = 270
$ws.Cells.Item(95, 11).Value = 17000
$ws.Cells.Item(95, 12).Value = 18000
$ws.Cells.Item(95, 13).Value = 17556
$ws.Cells.Item(95, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(95, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(95, 16).Value = 293
$ws.Cells.Item(95, 17).Value = 60
$ws.Cells.Item(95, 18).Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of column D.
$ws.Cells.Item(95, 4).NumberFormat = $ws.Cells.Item(96, 4).NumberFormat
